# Bug fix, Design and Exit button
# Append a new data row (row 10) to Sheet1 with the values used to flag
# the Design/Exit button behaviour: A10=1, B10=1, C10=1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(10, 1).Value = 1
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = 1
